# DbSchema NSKFDC-App.xlsx - Added table for reports upload
# Adds a new "UploadReports" table description in column A (rows 12-23),
# updates the scgjBatchNumber label to note it is unique, and tweaks
# the sheet view/row layout to match the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "UploadReports" schema table in column A ------------------------
# Row 12 is the table header (same styling as the other table headers,
# e.g. A5 "User", E12 "BatchDetails", G15 "CentreDetails", ...).
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Value = "UploadReports"

# Row 13 is the first field of the table (PK-style row styling, like A6).
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value = "reportId (Auto Increment)"

# Rows 14-21 are regular fields (styled like A7/A8/A9).
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = "reportType"

$ws.Range("A7").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Value = "reportPath"

$ws.Range("A7").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = "occupationCertificate  - Flag"

$ws.Range("A7").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = "attendanceSheet - Flag"

$ws.Range("A7").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = "nskfdcSheet - Flag"

$ws.Range("A7").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Value = "finalBatchReport - Flag"

$ws.Range("A7").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("A20").Value = "sdmsSheet - Flag"

$ws.Range("A7").Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Range("A21").Value = "selectionCommittee - Flag"

# Rows 22-23 are foreign-key references (styled like E9, the other FK row).
$ws.Range("E9").Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4122) | Out-Null
$ws.Range("A22").Value = "batchId - FK"

$ws.Range("E9").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value = "trainingPartnerEmail - FK"

$excel.CutCopyMode = 0

# --- scgjBatchNumber is now flagged as unique -----------------------------
$ws.Range("E21").Value = "scgjBatchNumber - unique"

# --- Cosmetic layout tweaks to match the authored workbook ---------------
# A short spacer row above the header row.
$ws.Rows.Item(3).RowHeight = 8.25

# View no longer frozen/scrolled to row 4, and selection moved to A7.
$ws.Range("A7").Select() | Out-Null
